$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Suite A (row 2) flips from Y to N
$ws.Range("C2").Value = "N"

# Suite E (row 6) flips from N to Y
$ws.Range("C6").Value = "Y"

# Suite F (row 7) flips from N to Y
$ws.Range("C7").Value = "Y"

# Update the active cell selection to C6
$ws.Range("C6").Select()
